$wb = $excel.ActiveWorkbook

# --- Rename 'variable' -> 'variable-code' and 'long_name' -> 'variable-label' ---
# (Close #144)

# "Codelists" sheet: header row uses variable / da_code_label / en_code_label / kl_code_label
$wsCodelists = $wb.Worksheets.Item("Codelists")
$wsCodelists.Range("D1").Value = "da_code-label"
$wsCodelists.Range("E1").Value = "en_code-label"
$wsCodelists.Range("F1").Value = "kl_code-label"

# "Variables" sheet: header row uses variable / da_long_name / en_long_name / kl_long_name
$wsVariables = $wb.Worksheets.Item("Variables")
$wsVariables.Range("C1").Value = "variable-code"

$wsCodelists.Range("A1").Value = "variable-code"

$wsVariables.Range("D1").Value = "da_variable-label"
$wsVariables.Range("E1").Value = "en_variable-label"
$wsVariables.Range("F1").Value = "kl_variable-label"

# --- Update view/selection state to match the saved workbook ---
$wsVariables.Activate() | Out-Null
$wsVariables.Range("D1").Select() | Out-Null

$wsCodelists.Activate() | Out-Null
$wsCodelists.Range("A2").Select() | Out-Null

# Window geometry (best effort; matches the on-disk <workbookView> values)
$excel.ActiveWindow.Left = 3810
$excel.ActiveWindow.Top = 3810
$excel.ActiveWindow.Width = 24210
$excel.ActiveWindow.Height = 21645
